$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Switch reference style to R1C1 (calcPr refMode="R1C1")
$excel.ReferenceStyle = 4  # xlR1C1

# New column F header + width
$ws.Columns.Item(6).ColumnWidth = 15

# Apply same style as column E to the new F column cells (rows 1-26)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E2:E26").Copy()
$ws.Range("F2:F26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 6).Value = "Sonstiges"

# "Ferien Experte" / "Ferien Betreuer" markers in column F
$ws.Cells.Item(5, 6).Value = "Ferien Betreuer"
$ws.Cells.Item(12, 6).Value = "Ferien Experte"
$ws.Cells.Item(16, 6).Value = "Ferien Experte"
$ws.Cells.Item(23, 6).Value = "Ferien Experte"

# Select F23 as active cell, matching diff's selection change
$ws.Range("F23").Select()
